$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style/format of the existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-19
$data = @(
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(9, 9),
    @(6, 6),
    @(4, 5),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

$excel.CutCopyMode = $false
